# Refresh the cryptos price/volume table (Sheet1, columns D "Price" and
# E "Volume(1h)") with the latest scraped figures.
#
# Some "Price" values are plain decimals (e.g. 304.47, 1.00, 18.77) even
# though the sheet always stores column D as text. A bare numeric-looking
# string assigned to .Value would be auto-coerced by Excel into a Number
# cell (and "1.00" would collapse to "1"), so those are entered with a
# leading apostrophe to force them to stay literal text, matching the
# original authoring. Values that already contain non-numeric punctuation
# (e.g. "42.978.70") or the Volume column's "  +0.08%  " strings are safe
# to set directly since Excel can't parse them as numbers anyway.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.978.70"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "2.305.34"
$ws.Range("E3").Value = "  +0.17%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'304.47"
$ws.Range("E5").Value = "  +1.23%  "
$ws.Range("D6").Value = "'97.23"
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("E7").Value = "  -2.00%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E9").Value = "  -1.12%  "
$ws.Range("D10").Value = "'35.38"
$ws.Range("E10").Value = "  -1.16%  "
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("D12").Value = "'18.77"
$ws.Range("E12").Value = "  +4.78%  "
$ws.Range("E13").Value = "  +1.57%  "
$ws.Range("E14").Value = "  +1.82%  "
$ws.Range("D15").Value = "2.665.16"
$ws.Range("E15").Value = "  +0.64%  "
$ws.Range("D16").Value = "2.308.09"
$ws.Range("E16").Value = "  +0.71%  "
$ws.Range("D17").Value = "'0.783"
$ws.Range("E17").Value = "  +0.70%  "
$ws.Range("D18").Value = "42.866.75"
$ws.Range("E18").Value = "  +0.06%  "
$ws.Range("E19").Value = "  -2.36%  "
$ws.Range("E20").Value = "  -0.91%  "
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("D22").Value = "'67.77"
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").Value = "'237.07"
$ws.Range("E23").Value = "  -1.69%  "
$ws.Range("D24").Value = "'2.16"
$ws.Range("E24").Value = "  +1.05%  "
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("D26").Value = "'2.43"
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").Value = "'24.89"
$ws.Range("E27").Value = "  -2.30%  "
$ws.Range("E28").Value = "  +7.69%  "
$ws.Range("D29").Value = "'165.99"
$ws.Range("E29").Value = "  +0.24%  "
$ws.Range("E30").Value = "  +0.15%  "
$ws.Range("D31").Value = "'32.82"
$ws.Range("E31").Value = "  -0.56%  "
$ws.Range("D33").Value = "'18.15"
$ws.Range("E33").Value = "  +6.38%  "
$ws.Range("E34").Value = "  -1.03%  "
$ws.Range("D35").Value = "'4.49"
$ws.Range("E35").Value = "  -7.94%  "
$ws.Range("E36").Value = "  -1.04%  "
$ws.Range("E37").Value = "  +0.46%  "
$ws.Range("E38").Value = "  -0.33%  "
$ws.Range("E39").Value = "  -0.25%  "
$ws.Range("E40").Value = "  +1.44%  "
$ws.Range("E41").Value = "  -0.59%  "
$ws.Range("D42").Value = "1.996.18"
$ws.Range("E42").Value = "  -1.19%  "
$ws.Range("E43").Value = "  -0.62%  "
$ws.Range("D44").Value = "'10.34"
$ws.Range("E44").Value = "  +1.96%  "
$ws.Range("E45").Value = "  +0.49%  "
$ws.Range("E46").Value = "  +3.51%  "
$ws.Range("E47").Value = "  -0.96%  "
$ws.Range("D48").Value = "2.531.26"
$ws.Range("E48").Value = "  +0.59%  "
$ws.Range("D49").Value = "'53.49"
$ws.Range("E49").Value = "  +0.05%  "
$ws.Range("E50").Value = "  -2.47%  "
$ws.Range("D51").Value = "'71.82"
$ws.Range("E51").Value = "  -0.38%  "
